$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "89.352.80"
$ws.Range("E2").Value2 = "  -1.43%  "

$ws.Range("D3").Value2 = "3.136.31"
$ws.Range("E3").Value2 = "  -1.33%  "

$ws.Range("E4").Value2 = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "214.78"
$ws.Range("E5").Value2 = "  +1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "635.62"
$ws.Range("E6").Value2 = "  +3.37%  "

$ws.Range("E7").Value2 = "  +1.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.796"
$ws.Range("E8").Value2 = "  +16.18%  "

$ws.Range("E9").Value2 = "  +0.09%  "

$ws.Range("D10").Value2 = "3.134.13"
$ws.Range("E10").Value2 = "  -0.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.563"
$ws.Range("E11").Value2 = "  -1.31%  "

$ws.Range("E12").Value2 = "  +1.50%  "

$ws.Range("E13").Value2 = "  -1.46%  "

$ws.Range("E14").Value2 = "  +2.89%  "

$ws.Range("D15").Value2 = "89.154.48"
$ws.Range("E15").Value2 = "  -1.23%  "

$ws.Range("D16").Value2 = "3.705.50"
$ws.Range("E16").Value2 = "  -1.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "32.22"
$ws.Range("E17").Value2 = "  -1.53%  "

$ws.Range("D18").Value2 = "3.150.40"
$ws.Range("E18").Value2 = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "3.43"
$ws.Range("E19").Value2 = "  +6.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "0.0000225"
$ws.Range("E20").Value2 = "  +18.60%  "

$ws.Range("E21").Value2 = "  -1.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "424.61"
$ws.Range("E22").Value2 = "  -2.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "8.43"
$ws.Range("E23").Value2 = "  -1.02%  "

$ws.Range("E24").Value2 = "  -3.63%  "

$ws.Range("E25").Value2 = "  +5.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "82.71"
$ws.Range("E26").Value2 = "  +10.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "11.51"
$ws.Range("E27").Value2 = "  -2.52%  "

$ws.Range("D28").Value2 = "3.297.92"
$ws.Range("E28").Value2 = "  -4.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.00"
$ws.Range("E29").Value2 = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.00"
$ws.Range("E30").Value2 = "  -0.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.157"
$ws.Range("E31").Value2 = "  -6.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.98"
$ws.Range("E32").Value2 = "  -3.96%  "

$ws.Range("E33").Value2 = "  -3.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "504.68"
$ws.Range("E34").Value2 = "  -5.48%  "

$ws.Range("E35").Value2 = "  +16.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "6.95"
$ws.Range("E36").Value2 = "  +0.98%  "

$ws.Range("E37").Value2 = "  +2.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "1.84"
$ws.Range("E38").Value2 = "  -1.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "22.31"
$ws.Range("E39").Value2 = "  +1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "22.29"
$ws.Range("E40").Value2 = "  -0.08%  "

$ws.Range("E41").Value2 = "  +0.45%  "

$ws.Range("E42").Value2 = "  +0.05%  "

$ws.Range("E43").Value2 = "  -2.28%  "

$ws.Range("E44").Value2 = "  -2.89%  "

$ws.Range("E45").Value2 = "  +8.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "146.03"
$ws.Range("E46").Value2 = "  +0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "43.73"
$ws.Range("E47").Value2 = "  -1.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "163.61"
$ws.Range("E48").Value2 = "  -5.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0652"
$ws.Range("E49").Value2 = "  +11.36%  "

$ws.Range("E50").Value2 = "  +3.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "24.15"
$ws.Range("E51").Value2 = "  +0.43%  "
